$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their literal text representation
# (e.g. trailing zeros / thousands-dot formatting) instead of being
# auto-coerced to numbers by Excel when assigned via .Value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.570.04'
$ws.Range("E2").Value = '  -1.97%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.884.31'
$ws.Range("E3").Value = '  -2.92%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.88'
$ws.Range("E5").Value = '  -2.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.21'
$ws.Range("E6").Value = '  +8.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.670'
$ws.Range("E7").Value = '  -2.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.752'
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("E10").Value = '  +5.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.41'
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.51'
$ws.Range("E13").Value = '  +2.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.512.02'
$ws.Range("E14").Value = '  -2.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.886.82'
$ws.Range("E15").Value = '  -2.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.95'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.97'
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("E18").Value = '  -3.99%  '
$ws.Range("E19").Value = '  -2.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.426.77'
$ws.Range("E20").Value = '  -1.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '441.28'
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("E22").Value = '  -6.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '94.20'
$ws.Range("E23").Value = '  -3.12%  '
$ws.Range("E24").Value = '  -3.78%  '
$ws.Range("E25").Value = '  -3.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.84'
$ws.Range("E26").Value = '  +4.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.03'
$ws.Range("E27").Value = '  -6.61%  '
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.53'
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.84'
$ws.Range("E30").Value = '  +14.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.19'
$ws.Range("E31").Value = '  -3.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.63'
$ws.Range("E32").Value = '  -2.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '48.38'
$ws.Range("E33").Value = '  +0.52%  '
$ws.Range("E34").Value = '  -3.75%  '
$ws.Range("E35").Value = '  +11.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '69.44'
$ws.Range("E36").Value = '  -3.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '634.13'
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.437'
$ws.Range("E38").Value = '  -0.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.146'
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.35'
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("E44").Value = '  -3.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.31'
$ws.Range("E46").Value = '  -3.76%  '
$ws.Range("E47").Value = '  -4.14%  '
$ws.Range("E48").Value = '  -14.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.910.88'
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("E50").Value = '  +3.74%  '
$ws.Range("E51").Value = '  -6.28%  '

# Row 43 and 45 swap (dogwifhat <-> Fetch.AI) plus value updates
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.17'
$ws.Range("E43").Value = '  +19.91%  '
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.86'
$ws.Range("E45").Value = '  +7.39%  '

Write-Output "done"